$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Timetable")

# Room renaming: old room names -> new room numbers
$ws.Range("C7").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("D7").Value = "Rehearsal with pianist`n(Room G22)"
$ws.Range("F7").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("B11").Value = "Private Lesson with Stephane RETY `n(Room G19)"
$ws.Range("E11").Value = "Private Lesson with Stephane RETY & pianist `n(Room G19)"

$ws.Range("D19").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("F19").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("B23").Value = "Ensemble `n(Room 245)"
$ws.Range("C23").Value = "Ensemble `n(Room 245)"
$ws.Range("D23").Value = "Ensemble `n(Room 245)"
$ws.Range("E23").Value = "Ensemble `n(Room 245)"
$ws.Range("F23").Value = "Ensemble `n(Room 245)"
